# Update gh-pages to output generated at 456a3b4
# Applies updated "want to go" counts, a renamed event title, and a
# refreshed cover image URL to both the "展览" and "全部类型" sheets.

$wb = $excel.ActiveWorkbook

$newTitle = "昆山·创世次元动漫游戏嘉年华之山海经兽兽盛会X梦幻启航X兽装盛宴邀您共赴！(免费展)"
$newCover = "//i1.hdslb.com/bfs/openplatform/202410/RpESHS911729246696391.jpeg"

# Map of worksheet name -> row number -> new F (想去人数) value
$sheetRows = @{
    "展览"   = @{ 4 = 1159; 8 = 271; 10 = 1027; 12 = 522; 13 = 551; 15 = 13063; 17 = 2; 19 = 5366; 20 = 5550; 21 = 4 }
    "全部类型" = @{ 4 = 1159; 24 = 271; 32 = 1027; 34 = 522; 35 = 551; 37 = 13063; 39 = 2; 42 = 5366; 43 = 5550; 44 = 4 }
}

# Rows (per sheet) whose event title (C) and cover image (I) also changed.
$titleRow = @{ "展览" = 10; "全部类型" = 32 }

foreach ($sheetName in $sheetRows.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)

    $rows = $sheetRows[$sheetName]
    foreach ($rowNum in $rows.Keys) {
        $ws.Range("F$rowNum").Value = $rows[$rowNum]
    }

    $trow = $titleRow[$sheetName]
    $ws.Range("C$trow").Value = $newTitle
    $ws.Range("I$trow").Value = $newCover
}

$wb.Save()
